$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("IRS-Cleared")
$ws1.Range("B2").Value = "ACUOSG8745"
$ws2 = $wb.Worksheets.Item("IRS-Bilateral")
$ws2.Range("B2").Value = "ACUOSG8745"
